# Update the FSA competition dates table on the active sheet ("Sheet1" in the
# Excel UI, which is the sheet carrying Table1 and the date/category rows).
#
# The edit swaps the dates for three groups of rows:
#   - Rows 65-71: 11 Sep 2022 (20220911) -> 18 Sep 2022 (20220918)
#   - Rows 72-75: 18 Sep 2022 (20220918) -> 11 Sep 2022 (20220911)
#   - Rows 76-81: 16 Oct 2022 (20221016) -> 8 Oct 2022 (20221008)
#
# Column K/L contain formulas that key off Table1[[#This Row],[Date]], so
# their cached values are refreshed automatically on recalculation once
# column A changes - no direct edits needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 65-71: 20220911 -> 20220918
$ws.Range("A65:A71").Value = 20220918

# Rows 72-75: 20220918 -> 20220911
$ws.Range("A72:A75").Value = 20220911

# Rows 76-81: 20221016 -> 20221008
$ws.Range("A76:A81").Value = 20221008

# Reflect the author's final cursor/scroll position in the saved view state.
$ws.Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$ws.Range("A81").Select() | Out-Null
